$wb = $excel.ActiveWorkbook

# "Overview" sheet: row 3 corresponds to the c558b3b0-... file.
# Its zh-cn (B3) and de-de (C3) status move from "In Translation" to
# "Ready for handoff", and the overall Latest Handoff Date (D3) is bumped.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-24 07:52:52"

# "zh-cn" sheet: row 3 (c558b3b0-... file) is now ready for handoff, and
# its Latest Handoff Datetime (E3) is stamped with the handoff time.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-24 07:52:48"

# "de-de" sheet: row 3 (c558b3b0-... file) is now ready for handoff, and
# its Latest Handoff Datetime (E3) is stamped with the handoff time.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-24 07:52:52"
